# Diet proportion up to predator age 20 & fully fixed max age bit
#
# Extends the UobsWtAge sheet with diet-proportion-by-weight rows for
# predator ages 16-20 (rows 77-101), mirroring the existing age-15 block
# (rows 72-76) and moves the sheet selection to E105.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UobsWtAge")

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
}

# --- Predator age 16 (rows 77-81) -----------------------------------
# Exact repeat of the age-15 block (rows 72-76) above it, just with
# Pred_age (column E) bumped from 15 to 16. Same stomach-proportion
# values, including the scientific-notation formatted Prey_age=2 row.
Set-Row 77 1 1 0 0 16 1 0 10 0.111619969468499
Set-Row 78 1 1 0 0 16 2 0 10 0.0000099994000359978393
Set-Row 79 1 1 0 0 16 3 0 10 0.0296082235065896
Set-Row 80 1 1 0 0 16 4 0 10 0.0308081515109093
Set-Row 81 1 1 0 0 16 5 0 10 0.0281083135011899

# I78 keeps the scientific-notation number format used by the same
# column on row 73 (Prey_age = 2 rows render with 0.00E+00).
$ws.Range("I78").NumberFormat = "0.00E+00"

# --- Predator ages 17-20 (rows 82-101) -------------------------------
# Same five stomach-proportion values (rounded) repeated for each
# predator age, rendered with the new plain-black Calibri font.
$ages = 17, 18, 19, 20
$r = 82
foreach ($age in $ages) {
    Set-Row ($r + 0) 1 1 0 0 $age 1 0 10 0.11161997
    Set-Row ($r + 1) 1 1 0 0 $age 2 0 10 0.00001
    Set-Row ($r + 2) 1 1 0 0 $age 3 0 10 0.029608220000000001
    Set-Row ($r + 3) 1 1 0 0 $age 4 0 10 0.030808149999999999
    Set-Row ($r + 4) 1 1 0 0 $age 5 0 10 0.028108310000000001
    $r = $r + 5
}

# Apply the new plain-black-Calibri font to the whole A82:I101 block.
$ws.Range("A82:I101").Font.Color = 0

# The Prey_age = 2 rows (83, 88, 93, 98) additionally use scientific
# notation in column I.
$ws.Range("I83").NumberFormat = "0.00E+00"
$ws.Range("I88").NumberFormat = "0.00E+00"
$ws.Range("I93").NumberFormat = "0.00E+00"
$ws.Range("I98").NumberFormat = "0.00E+00"

# --- Sheet selection --------------------------------------------------
[void]$ws.Range("E105").Select()
